# Deploy the implementation guide.
# Updates the Metadata sheet: refresh the generation Date, replace the
# Contact value, and insert a new "Jurisdiction" property row right after
# "Contact" (pushing Description..Count down by one row).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# 1. Refresh the "Date" value (row 8, column B).
$ws.Range("B8").Value = "2024-10-02T15:04:17+00:00"

# 2. Replace the "Contact" value (row 10, column B).
$ws.Range("B10").Value = "Ferlab.bio (http://example.org/example-publisher)"

# 3. Make room for a new row right below "Contact" by inserting a blank row
#    at 11, shifting the existing rows 11:21 down to 12:22.
$ws.Range("A11:B11").Insert(-4121)

# 4. The inserted row doesn't inherit the table's row styling, so copy the
#    formatting (borders/fill/alignment) from the row directly below.
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 5. Populate the freed-up row 11 with the new "Jurisdiction" property (no
#    value supplied).
$ws.Range("A11").Value = "Jurisdiction"
$ws.Range("B11").Value = ""
